$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows the data columns (Fecha, Volumen, Precio minimo/maximo/
# promedio ponderado, Origen, Precio $/Kg) being redistributed among rows
# 2, 4, 5, 7, 8, 9, 10, 11 -- the other columns (Mercado ID, Mercado,
# Region, Codreg, Tipo, Producto ID, Producto, Categoria ID, Categoria,
# Variedad, Calidad, Unidad de comercializacion, Kg / unidad) stay as-is.
#
# Target values per row for columns D (Fecha), M (Volumen), N (Precio
# minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen), S (Precio $/Kg):

$updates = @{
    2  = @{ D = 44574; M = 200; N = 3000; O = 3000; P = 3000; R = "Región de La Araucanía"; S = 3000 }
    4  = @{ D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; R = "Región de La Araucanía"; S = 2800 }
    5  = @{ D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó";    S = 5000 }
    7  = @{ D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; R = "Región de La Araucanía"; S = 2400 }
    8  = @{ D = 44551; M = 120; N = 4500; O = 4500; P = 4500; R = "Región de O'Higgins";    S = 4500 }
    9  = @{ D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; R = "Región de O'Higgins";    S = 3000 }
    10 = @{ D = 44592; M = 5;   N = 7500; O = 7500; P = 7500; R = "Región de La Araucanía"; S = 7500 }
    11 = @{ D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; R = "Región de La Araucanía"; S = 1800 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
